$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 3 (Doru1): Role -> Professor, Title -> Eng., Department -> AC
$ws.Range("C3").Value = "Professor"
$ws.Range("I3").Value = "Eng."
$ws.Range("H3").Value = "AC"

# Row 5 (Doru3): Role -> Campus_Student
$ws.Range("C5").Value = "Campus_Student"

# Update active selection to C5
$ws.Range("C5").Select()
